$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2127.2
$ws.Range("I38").Value = 47.42857
$ws.Range("J38").Value = 6980
$ws.Range("K38").Value = 142.28571
$ws.Range("L38").Value = 20940
$ws.Range("M38").Value = 229.71429
$ws.Range("N38").Value = -21684
$ws.Range("H40").Value = 7372.909
$ws.Range("I40").Value = 5872.2856
$ws.Range("K40").Value = 5872.2856
$ws.Range("M40").Value = -5697.2856
$ws.Range("H41").Value = 964.8889
$ws.Range("I41").Value = 835.5
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 835.5
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = -395.5
$ws.Range("N41").Value = -2880
$ws.Range("H98").Value = 817.25
$ws.Range("I98").Value = 689.8333
$ws.Range("K98").Value = 689.8333
$ws.Range("M98").Value = 808.1667
$ws.Range("H122").Value = 817.25
$ws.Range("I122").Value = 689.8333
$ws.Range("K122").Value = 2069.4999
$ws.Range("M122").Value = 380.5001000000002
$ws.Range("H137").Value = 1530.6875
$ws.Range("I137").Value = 771.44446
$ws.Range("J137").Value = 2506.8572
$ws.Range("K137").Value = 2314.33338
$ws.Range("L137").Value = 7520.571599999999
$ws.Range("M137").Value = 235.66662
$ws.Range("N137").Value = -12620.5716
$ws.Range("H138").Value = 2717.5715
$ws.Range("I138").Value = 1555
$ws.Range("K138").Value = 4665
$ws.Range("M138").Value = 475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1609.9333
$ws.Range("I61").Value = 1609.9333
$ws.Range("K61").Value = 1609.9333
$ws.Range("M61").Value = -1397.9333
$ws.Range("H122").Value = 2646
$ws.Range("I122").Value = 2646
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7938
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5488
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 1609.9333
$ws.Range("I136").Value = 1609.9333
$ws.Range("K136").Value = 4829.7999
$ws.Range("M136").Value = -2279.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1733
$ws.Range("I20").Value = 669.3333
$ws.Range("J20").Value = 2371.2
$ws.Range("K20").Value = 669.3333
$ws.Range("L20").Value = 2371.2
$ws.Range("M20").Value = -422.3333
$ws.Range("N20").Value = -2865.2
$ws.Range("H107").Value = 55564400
$ws.Range("I107").Value = 166673260
$ws.Range("K107").Value = 166673260
$ws.Range("M107").Value = -166671340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.40000000000001
$ws.Range("I7").Value = 70.5
$ws.Range("J7").Value = 67.75
$ws.Range("K7").Value = 70.5
$ws.Range("L7").Value = 67.75
$ws.Range("M7").Value = 42.5
$ws.Range("N7").Value = -293.75
$ws.Range("H31").Value = 1211.0555
$ws.Range("J31").Value = 2749.5
$ws.Range("L31").Value = 2749.5
$ws.Range("N31").Value = -3339.5
$ws.Range("H34").Value = 1211.0555
$ws.Range("J34").Value = 2749.5
$ws.Range("L34").Value = 2749.5
$ws.Range("N34").Value = -3153.5
$ws.Range("H62").Value = 2511.5
$ws.Range("I62").Value = 799
$ws.Range("J62").Value = 3082.3333
$ws.Range("K62").Value = 799
$ws.Range("L62").Value = 3082.3333
$ws.Range("M62").Value = -175
$ws.Range("N62").Value = -4330.3333
$ws.Range("H65").Value = 2511.5
$ws.Range("I65").Value = 799
$ws.Range("J65").Value = 3082.3333
$ws.Range("K65").Value = 3995
$ws.Range("L65").Value = 15411.6665
$ws.Range("M65").Value = -875
$ws.Range("N65").Value = -21651.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 138.66667
$ws.Range("I12").Value = 22.75
$ws.Range("K12").Value = 68.25
$ws.Range("M12").Value = 104.75
$ws.Range("H98").Value = 119.833336
$ws.Range("I98").Value = 127.5
$ws.Range("K98").Value = 382.5
$ws.Range("M98").Value = 1115.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 967.7273
$ws.Range("I16").Value = 967.7273
$ws.Range("K16").Value = 967.7273
$ws.Range("M16").Value = -797.7273
$ws.Range("H22").Value = 740.0909
$ws.Range("I22").Value = 230.33333
$ws.Range("J22").Value = 931.25
$ws.Range("K22").Value = 230.33333
$ws.Range("L22").Value = 931.25
$ws.Range("M22").Value = 64.66667000000001
$ws.Range("N22").Value = -1521.25
$ws.Range("H27").Value = 740.0909
$ws.Range("I27").Value = 230.33333
$ws.Range("J27").Value = 931.25
$ws.Range("K27").Value = 230.33333
$ws.Range("L27").Value = 931.25
$ws.Range("M27").Value = -123.33333
$ws.Range("N27").Value = -1145.25
$ws.Range("H55").Value = 1196.8
$ws.Range("I55").Value = 1331
$ws.Range("J55").Value = 1079.375
$ws.Range("K55").Value = 1331
$ws.Range("L55").Value = 1079.375
$ws.Range("M55").Value = -1158
$ws.Range("N55").Value = -1425.375
$ws.Range("H68").Value = 2239
$ws.Range("I68").Value = 1298.875
$ws.Range("J68").Value = 5999.5
$ws.Range("K68").Value = 1298.875
$ws.Range("L68").Value = 5999.5
$ws.Range("M68").Value = -549.875
$ws.Range("N68").Value = -7497.5
$ws.Range("H71").Value = 2239
$ws.Range("I71").Value = 1298.875
$ws.Range("J71").Value = 5999.5
$ws.Range("K71").Value = 6494.375
$ws.Range("L71").Value = 29997.5
$ws.Range("M71").Value = -2750.375
$ws.Range("N71").Value = -37485.5
$ws.Range("H93").Value = 1534.0667
$ws.Range("I93").Value = 1858.7142
$ws.Range("J93").Value = 1250
$ws.Range("K93").Value = 1858.7142
$ws.Range("L93").Value = 1250
$ws.Range("M93").Value = -610.7141999999999
$ws.Range("N93").Value = -3746
$ws.Range("H132").Value = 4799.857
$ws.Range("I132").Value = 4799.857
$ws.Range("K132").Value = 14399.571
$ws.Range("M132").Value = -11869.571
$ws.Range("H136").Value = 2999.5
$ws.Range("I136").Value = 2999.5
$ws.Range("K136").Value = 8998.5
$ws.Range("M136").Value = -6448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H126").Value = 5855
$ws.Range("I126").Value = 492.5
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 1477.5
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = 992.5
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 2999
$ws.Range("I132").Value = 2999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8997
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6467
$ws.Range("N132").ClearContents()
